{"js": "// Apply a series of small wording fixes throughout the requirements\n// document. Each fix is performed with a targeted search-and-replace so\n// that unrelated text (including other runs/paragraphs with similar\n// wording) is left untouched.\n\nasync function replaceOnce(context, searchText, replacementText, options) {\n  const body = context.document.body;\n  const results = body.search(searchText, Object.assign({ matchCase: true, matchWholeWord: false }, options || {}));\n  results.load(\"items\");\n  await context.sync();\n\n  if (results.items.length !== 1) {\n    throw new Error(\n      \"Expected exactly 1 match for \" + JSON.stringify(searchText) + \" but found \" + results.items.length\n    );\n  }\n\n  results.items[0].insertText(replacementText, Word.InsertLocation.replace);\n  await context.sync();\n}\n\n// 1. \"...log unencrypted packets and where the packets were sent to.\" ->\n//    \"...log unencrypted packets and where the packets were sent.\"\nawait replaceOnce(\n  context,\n  \"unencrypted packets and where the packets were sent to\",\n  \"unencrypted packets and where the packets were sent\"\n);\n\n// 2. \"Packets should be checked and logged within 110 milliseconds.\" ->\n//    \"...within 100 milliseconds.\"\nawait replaceOnce(\n  context,\n  \"Packets should be checked and logged within 110 milliseconds\",\n  \"Packets should be checked and logged within 100 milliseconds\"\n);\n\n// 3. \"Fit Criterion: Administrators should be able to configure and\n//    monitor the analysis module with minimal training.\" ->\n//    \"Fit Criterion: Users should be able to configure and monitor the\n//    analysis module with minimal knowledge.\"\nawait replaceOnce(\n  context,\n  \"Fit Criterion: Administrators should be able to configure and monitor the analysis module with minimal training.\",\n  \"Fit Criterion: Users should be able to configure and monitor the analysis module with minimal knowledge.\"\n);\n\n// 4. \"Rationale: Detailed reports aid in post-incident analysis and audit\n//    trails.\" -> \"Rationale: Detailed reports aid in the analysis after\n//    the fact.\"\nawait replaceOnce(\n  context,\n  \"Rationale: Detailed reports aid in post-incident analysis and audit trails.\",\n  \"Rationale: Detailed reports aid in the analysis after the fact.\"\n);\n\n// 5. \"...apply them automatically with administrator approval.\" ->\n//    \"...apply them automatically with user approval.\"\nawait replaceOnce(\n  context,\n  \"Fit Criterion: The system should check for updates periodically and apply them automatically with administrator approval.\",\n  \"Fit Criterion: The system should check for updates periodically and apply them automatically with user approval.\"\n);\n\n// 6. \"Description: The system should be able to provide real-time\n//    monitoring...\" -> \"Description: The system should provide\n//    real-time monitoring...\"\nawait replaceOnce(\n  context,\n  \"Description: The system should be able to provide real-time monitoring and alerting capabilities with minimal latency.\",\n  \"Description: The system should provide real-time monitoring and alerting capabilities with minimal latency.\"\n);\n\n// 7. \"...the code should be easy to read so that everyone can understand\n//    what each part does.\" -> \"...so everyone can understand...\"\nawait replaceOnce(\n  context,\n  \"the code should be easy to read so that everyone can understand what\",\n  \"the code should be easy to read so everyone can understand what\"\n);\n", "ps1": "# Apply a series of small wording fixes throughout the requirements\n# document using Find/Replace on the document's main story range.\n\n$wdFindContinue = 1\n$wdReplaceAll = 2\n\n$d = $word.ActiveDocument\n\nfunction Replace-Text($findText, $replaceText) {\n    $find = $d.Content.Find\n    $find.ClearFormatting()\n    $find.Replacement.ClearFormatting()\n    $find.Execute(\n        $findText,    # FindText\n        $true,        # MatchCase\n        $false,       # MatchWholeWord\n        $false,       # MatchWildcards\n        $false,       # MatchSoundsLike\n        $false,       # MatchAllWordForms\n        $true,        # Forward\n        $wdFindContinue, # Wrap\n        $false,       # Format\n        $replaceText, # ReplaceWith\n        $wdReplaceAll # Replace\n    ) | Out-Null\n}\n\n# 1. \"...log unencrypted packets and where the packets were sent to.\" ->\n#    \"...log unencrypted packets and where the packets were sent.\"\nReplace-Text \"unencrypted packets and where the packets were sent to\" \"unencrypted packets and where the packets were sent\"\n\n# 2. \"Packets should be checked and logged within 110 milliseconds.\" ->\n#    \"...within 100 milliseconds.\"\nReplace-Text \"Packets should be checked and logged within 110 milliseconds\" \"Packets should be checked and logged within 100 milliseconds\"\n\n# 3. \"Fit Criterion: Administrators should be able to configure and\n#    monitor the analysis module with minimal training.\" ->\n#    \"Fit Criterion: Users should be able to configure and monitor the\n#    analysis module with minimal knowledge.\"\nReplace-Text \"Fit Criterion: Administrators should be able to configure and monitor the analysis module with minimal training.\" \"Fit Criterion: Users should be able to configure and monitor the analysis module with minimal knowledge.\"\n\n# 4. \"Rationale: Detailed reports aid in post-incident analysis and audit\n#    trails.\" -> \"Rationale: Detailed reports aid in the analysis after\n#    the fact.\"\nReplace-Text \"Rationale: Detailed reports aid in post-incident analysis and audit trails.\" \"Rationale: Detailed reports aid in the analysis after the fact.\"\n\n# 5. \"...apply them automatically with administrator approval.\" ->\n#    \"...apply them automatically with user approval.\"\nReplace-Text \"Fit Criterion: The system should check for updates periodically and apply them automatically with administrator approval.\" \"Fit Criterion: The system should check for updates periodically and apply them automatically with user approval.\"\n\n# 6. \"Description: The system should be able to provide real-time\n#    monitoring...\" -> \"Description: The system should provide\n#    real-time monitoring...\"\nReplace-Text \"Description: The system should be able to provide real-time monitoring and alerting capabilities with minimal latency.\" \"Description: The system should provide real-time monitoring and alerting capabilities with minimal latency.\"\n\n# 7. \"...the code should be easy to read so that everyone can understand\n#    what each part does.\" -> \"...so everyone can understand...\"\nReplace-Text \"the code should be easy to read so that everyone can understand what\" \"the code should be easy to read so everyone can understand what\"\n"}
